{"js": "// Update the date line (first paragraph of the document body) and the\n// twenty two-digit-by-two-digit multiplication problems living in the\n// first table's non-blank rows (0, 4, 9, 14, 19).\n\nconst body = context.document.body;\n\n// --- 1) Update the date paragraph (\"2025-04-11 Friday\" -> \"2025-04-12 Saturday\") ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst datePara = paragraphs.items[0];\ndatePara.insertText(\"2025-04-12 Saturday\", Word.InsertLocation.replace);\n\n// --- 2) Update the multiplication problems inside the first table ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// (rowIndex, columnIndex, newValue) for every populated cell, in document order.\nconst updates = [\n  [0, 0, \"31\u00d788=\"],\n  [0, 1, \"24\u00d775=\"],\n  [0, 2, \"91\u00d771=\"],\n  [0, 3, \"15\u00d787=\"],\n  [0, 4, \"96\u00d719=\"],\n  [4, 0, \"64\u00d758=\"],\n  [4, 1, \"53\u00d791=\"],\n  [4, 2, \"58\u00d731=\"],\n  [4, 3, \"81\u00d796=\"],\n  [4, 4, \"54\u00d753=\"],\n  [9, 0, \"65\u00d795=\"],\n  [9, 1, \"51\u00d761=\"],\n  [9, 2, \"55\u00d737=\"],\n  [9, 3, \"96\u00d790=\"],\n  [9, 4, \"85\u00d751=\"],\n  [14, 0, \"43\u00d745=\"],\n  [14, 1, \"38\u00d757=\"],\n  [14, 2, \"50\u00d785=\"],\n  [14, 3, \"92\u00d754=\"],\n  [14, 4, \"54\u00d775=\"],\n  [19, 0, \"41\u00d746=\"],\n  [19, 1, \"51\u00d766=\"],\n  [19, 2, \"34\u00d785=\"],\n  [19, 3, \"16\u00d743=\"],\n  [19, 4, \"40\u00d712=\"],\n];\n\nfor (const [rowIndex, columnIndex, newValue] of updates) {\n  const cell = table.getCell(rowIndex, columnIndex);\n  cell.value = newValue;\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the twenty two-digit-by-two-digit multiplication\n# problems living in the first table's non-blank rows (1, 5, 10, 15, 20 -\n# Word COM uses 1-based row/column indices).\n\n$d = $word.ActiveDocument\n\n# --- 1) Update the date paragraph (\"2025-04-11 Friday\" -> \"2025-04-12 Saturday\") ---\n$p = $d.Paragraphs.Item(1)\n$p.Range.Text = \"2025-04-12 Saturday\"\n\n# --- 2) Update the multiplication problems inside the first table ---\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text = \"31\u00d788=\"\n$t.Cell(1, 2).Range.Text = \"24\u00d775=\"\n$t.Cell(1, 3).Range.Text = \"91\u00d771=\"\n$t.Cell(1, 4).Range.Text = \"15\u00d787=\"\n$t.Cell(1, 5).Range.Text = \"96\u00d719=\"\n\n$t.Cell(5, 1).Range.Text = \"64\u00d758=\"\n$t.Cell(5, 2).Range.Text = \"53\u00d791=\"\n$t.Cell(5, 3).Range.Text = \"58\u00d731=\"\n$t.Cell(5, 4).Range.Text = \"81\u00d796=\"\n$t.Cell(5, 5).Range.Text = \"54\u00d753=\"\n\n$t.Cell(10, 1).Range.Text = \"65\u00d795=\"\n$t.Cell(10, 2).Range.Text = \"51\u00d761=\"\n$t.Cell(10, 3).Range.Text = \"55\u00d737=\"\n$t.Cell(10, 4).Range.Text = \"96\u00d790=\"\n$t.Cell(10, 5).Range.Text = \"85\u00d751=\"\n\n$t.Cell(15, 1).Range.Text = \"43\u00d745=\"\n$t.Cell(15, 2).Range.Text = \"38\u00d757=\"\n$t.Cell(15, 3).Range.Text = \"50\u00d785=\"\n$t.Cell(15, 4).Range.Text = \"92\u00d754=\"\n$t.Cell(15, 5).Range.Text = \"54\u00d775=\"\n\n$t.Cell(20, 1).Range.Text = \"41\u00d746=\"\n$t.Cell(20, 2).Range.Text = \"51\u00d766=\"\n$t.Cell(20, 3).Range.Text = \"34\u00d785=\"\n$t.Cell(20, 4).Range.Text = \"16\u00d743=\"\n$t.Cell(20, 5).Range.Text = \"40\u00d712=\"\n"}
